$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = '(''Celestial Colonnade'', [''Land'', ''Celestial Colonnade enters the battlefield tapped.'', ''{T}: Add {W} or {U}.'', ''{3}{W}{U}: Until end of turn, Celestial Colonnade becomes a 4/4 white and blue Elemental creature with flying and vigilance. It’s still a land.''])'
$ws.Range("A3").Value = '(''Comet Storm'', [''{X}{R}{R}'', ''Instant'', ''Multikicker {1} (You may pay an additional {1} any number of times as you cast this spell.)'', ''Choose any target, then choose another target for each time this spell was kicked. Comet Storm deals X damage to each of them.''])'
$ws.Range("A4").Value = '(''Hada Freeblade'', [''{W}'', ''Creature — Human Soldier Ally'', ''Whenever Hada Freeblade or another Ally enters the battlefield under your control, you may put a +1/+1 counter on Hada Freeblade.'', ''0/1''])'
$ws.Range("A5").Value = '(''Joraga Warcaller'', [''{G}'', ''Creature — Elf Warrior'', ''Multikicker {1}{G} (You may pay an additional {1}{G} any number of times as you cast this spell.)'', ''Joraga Warcaller enters the battlefield with a +1/+1 counter on it for each time it was kicked.'', ''Other Elf creatures you control get +1/+1 for each +1/+1 counter on Joraga Warcaller.'', ''1/1''])'
$ws.Range("A6").Value = '(''Kalastria Highborn'', [''{B}{B}'', ''Creature — Vampire Shaman'', ''Whenever Kalastria Highborn or another Vampire you control dies, you may pay {B}. If you do, target player loses 2 life and you gain 2 life.'', ''2/2''])'
$ws.Range("A7").Value = '(''Ruthless Cullblade'', [''{1}{B}'', ''Creature — Vampire Warrior'', ''Ruthless Cullblade gets +2/+1 as long as an opponent has 10 or less life.'', ''2/1''])'

$ws.Rows("8:33").Delete()
